$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("year")
$ws.Activate()

# Delete the two rows that are being dropped from the scaling-year table:
#   row 7: rou / all / E=2010
#   row 8: mkd / all / E=1980
$ws.Rows("7:8").Delete()

# Row 2 becomes a generic "all/all" row limiting scaling to 1992-2009,
# replacing the old idn/1A2 special case, and drops its comment.
$ws.Range("A2").Value = "all"
$ws.Range("B2").Value = "all"
$ws.Range("F2").Value = 1992
$ws.Range("G2").Value = 2009
$ws.Range("H2").ClearContents()

# Row 3 (svk): only the end year changes.
$ws.Range("G3").Value = 2009

# Rows 4-6 (rou): shift the selected scaling years forward.
$ws.Range("E4").Value = 1992
$ws.Range("E5").Value = 2000
$ws.Range("E6").Value = 2010

# Rows 8-30: FSU + former-Yugoslavia countries now stop scaling at 2009
# instead of 2020.
$ws.Range("G8:G30").Value = 2009

# Row 31 (rail, all sectors 1A3c): narrow the scaling window.
$ws.Range("F31").Value = 2000
$ws.Range("G31").Value = 2009

# Row 32 (rail, alb) is unchanged in value, nothing to update.

# Restore the view: frozen header row, selection on G34 (last data row).
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G34").Select()
